$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata" ---
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/ValueSet/eng-product"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"
$meta.Range("B11").Value = "LinuxForHealth engagement products"

# --- Sheet "Include from Engagement Produ" ---
$inc = $wb.Worksheets.Item("Include from Engagement Produ")

$inc.Range("B4").Value = "http://linuxforhealth.org/fhir/cdm/CodeSystem/eng-product"
